$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "MODEL_CONDITION" header text to "MODELCONDITION".
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# 2. Column A in the original sheet only duplicated the GENE column (0/8/10)
#    and carried the header style with no header text. Deleting it shifts
#    columns B:F left into A:E, producing the new A1:E4 layout where
#    EL_Astral_exact100 / FNRATE_ASTRAL / TAXON / MODELCONDITION / GENE are
#    the headers in A1:E1.
$ws.Columns("A").Delete()
